$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.30"
$ws.Range("E2").Value = "'-0.75%"
$ws.Range("G2").Value = "'23"

$ws.Range("D3").Value = "'36.58"
$ws.Range("E3").Value = "'2.68%"
$ws.Range("G3").Value = "'23"

$ws.Range("D4").Value = "'4.982"
$ws.Range("E4").Value = "'-2.02%"
$ws.Range("G4").Value = "'23"

$ws.Range("D5").Value = "'0.07699"
$ws.Range("E5").Value = "'-0.65%"
$ws.Range("G5").Value = "'23"

$ws.Range("D6").Value = "'2.080"
$ws.Range("E6").Value = "'-7.09%"
$ws.Range("G6").Value = "'23"

$ws.Range("D7").Value = "'7.914"
$ws.Range("E7").Value = "'-1.56%"
$ws.Range("G7").Value = "'23"

$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.022"
$ws.Range("E8").Value = "'-0.69%"
$ws.Range("G8").Value = "'23"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9202"
$ws.Range("E9").Value = "'-0.79%"
$ws.Range("G9").Value = "'23"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09719"
$ws.Range("E10").Value = "'2.79%"
$ws.Range("G10").Value = "'23"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1858"
$ws.Range("E11").Value = "'1.29%"
$ws.Range("G11").Value = "'23"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08526"
$ws.Range("E12").Value = "'0.12%"
$ws.Range("G12").Value = "'23"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03516"
$ws.Range("E13").Value = "'-3.71%"
$ws.Range("G13").Value = "'23"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09958"
$ws.Range("E14").Value = "'0.31%"
$ws.Range("G14").Value = "'23"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001461"
$ws.Range("E15").Value = "'-1.83%"
$ws.Range("G15").Value = "'23"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005624"
$ws.Range("E16").Value = "'-2.19%"
$ws.Range("G16").Value = "'23"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.74%"
$ws.Range("G17").Value = "'23"

$ws.Range("D18").Value = "'2.406"
$ws.Range("E18").Value = "'10.16%"
$ws.Range("G18").Value = "'23"

$ws.Range("D19").Value = "'0.3387"
$ws.Range("E19").Value = "'-2.19%"
$ws.Range("G19").Value = "'23"

$ws.Range("D20").Value = "'0.1333"
$ws.Range("E20").Value = "'0.73%"
$ws.Range("G20").Value = "'23"

$ws.Range("D21").Value = "'4.779"
$ws.Range("E21").Value = "'4.77%"
$ws.Range("G21").Value = "'23"

$ws.Range("D22").Value = "'0.2196"
$ws.Range("E22").Value = "'-1.87%"
$ws.Range("G22").Value = "'23"

$ws.Range("D23").Value = "'0.04586"
$ws.Range("E23").Value = "'-1.85%"
$ws.Range("G23").Value = "'23"

$ws.Range("D24").Value = "'0.005088"
$ws.Range("E24").Value = "'12.54%"
$ws.Range("G24").Value = "'23"

$ws.Range("D25").Value = "'0.001228"
$ws.Range("E25").Value = "'-0.60%"
$ws.Range("G25").Value = "'23"

$ws.Range("D26").Value = "'0.0001394"
$ws.Range("E26").Value = "'6.77%"
$ws.Range("G26").Value = "'23"

$ws.Range("G27").Value = "'23"

$ws.Range("G28").Value = "'23"

$ws.Range("G29").Value = "'23"

$ws.Range("G30").Value = "'23"

$ws.Range("G31").Value = "'23"

$ws.Range("G32").Value = "'23"

$ws.Range("G33").Value = "'23"

$ws.Range("G34").Value = "'23"

$ws.Range("G35").Value = "'23"

$ws.Range("G36").Value = "'23"

$ws.Range("G37").Value = "'23"

$ws.Range("G38").Value = "'23"

$ws.Range("D39").Value = "'0.01762"
$ws.Range("E39").Value = "'-0.06%"
$ws.Range("G39").Value = "'23"

$ws.Range("D40").Value = "'0.04643"
$ws.Range("E40").Value = "'-1.80%"
$ws.Range("G40").Value = "'23"

$ws.Range("D41").Value = "'0.007418"
$ws.Range("E41").Value = "'-6.74%"
$ws.Range("G41").Value = "'23"

$ws.Range("D42").Value = "'0.1391"
$ws.Range("E42").Value = "'-1.22%"
$ws.Range("G42").Value = "'23"

$ws.Range("D43").Value = "'0.007706"
$ws.Range("E43").Value = "'-2.14%"
$ws.Range("G43").Value = "'23"

$ws.Range("D44").Value = "'0.002240"
$ws.Range("E44").Value = "'0.43%"
$ws.Range("G44").Value = "'23"

$ws.Range("D45").Value = "'0.01030"
$ws.Range("E45").Value = "'6.77%"
$ws.Range("G45").Value = "'23"

$ws.Range("D46").Value = "'0.00006274"
$ws.Range("E46").Value = "'0.97%"
$ws.Range("G46").Value = "'23"

$ws.Range("E47").Value = "'-0.62%"
$ws.Range("G47").Value = "'23"

$ws.Range("D48").Value = "'0.0005787"
$ws.Range("E48").Value = "'-0.23%"
$ws.Range("G48").Value = "'23"

$ws.Range("D49").Value = "'35.32"
$ws.Range("E49").Value = "'510.35%"
$ws.Range("G49").Value = "'23"

$ws.Range("D50").Value = "'0.001996"
$ws.Range("E50").Value = "'-26.10%"
$ws.Range("G50").Value = "'23"

$ws.Range("D51").Value = "'0.00002096"
$ws.Range("E51").Value = "'-0.62%"
$ws.Range("G51").Value = "'23"
